$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 5689
$ws.Range("F9").Value = 534
$ws.Range("F11").Value = 1546
$ws.Range("F12").Value = 9
$ws.Range("F14").Value = 649
$ws.Range("F15").Value = 1537
$ws.Range("F16").Value = 1537
$ws.Range("F17").Value = 1423
$ws.Range("F18").Value = 327
$ws.Range("F19").Value = 35
$ws.Range("F20").Value = 561
$ws.Range("F21").Value = 4069
$ws.Range("F22").Value = 4069
$ws.Range("F23").Value = 663
$ws.Range("F25").Value = 776
$ws.Range("F26").Value = 36
$ws.Range("F27").Value = 2235
$ws.Range("F28").Value = 35
$ws.Range("F29").Value = 317
$ws.Range("F31").Value = 39
$ws.Range("F32").Value = 1182
$ws.Range("F33").Value = 773
$ws.Range("F35").Value = 1090
$ws.Range("F36").Value = 1100

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 99
$ws.Range("F18").Value = 275
$ws.Range("F20").Value = 480

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 556
$ws.Range("F5").Value = 96
$ws.Range("F6").Value = 188

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F8").Value = 556
$ws.Range("F9").Value = 96
$ws.Range("F10").Value = 5689
$ws.Range("F16").Value = 99
$ws.Range("F18").Value = 534
$ws.Range("F21").Value = 1546
$ws.Range("F23").Value = 9
$ws.Range("F25").Value = 1537
$ws.Range("F27").Value = 1423
$ws.Range("F28").Value = 327
$ws.Range("F29").Value = 35
$ws.Range("F30").Value = 561
$ws.Range("F32").Value = 4069
$ws.Range("F33").Value = 4069
$ws.Range("F34").Value = 663
$ws.Range("F36").Value = 776
$ws.Range("F37").Value = 36
$ws.Range("F38").Value = 2235
$ws.Range("F39").Value = 35
$ws.Range("F41").Value = 39
$ws.Range("F44").Value = 275
$ws.Range("F46").Value = 480
$ws.Range("F47").Value = 773
$ws.Range("F49").Value = 1090
$ws.Range("F50").Value = 1100

